$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Project Planner")

# --- Update project title (B1) ---
$ws.Range("B1").Value = "Project 20 Planner"

# --- Update selected period (H2) ---
$ws.Range("H2").Value = 30

# --- Fill in Actual Start / Actual Duration / % Complete figures ---
# Row 5: 1. Research (section header) - only Actual Start filled
$ws.Range("E5").Value = 1

# Row 6: 1.1 Background research
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 4
$ws.Range("G6").Value = 1

# Row 7: 1.2 Group page
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 2
$ws.Range("G7").Value = 1

# Row 8: 1.3 Technical research
$ws.Range("E8").Value = 7
$ws.Range("G8").Value = 0.25

# Row 9: 1.3 User research
$ws.Range("E9").Value = 6
$ws.Range("F9").Value = 2
$ws.Range("G9").Value = 1

# Row 10: 2. Design (section header) - only Actual Start filled
$ws.Range("E10").Value = 6

# Row 11: 2.1 Write requirements
$ws.Range("E11").Value = 7
$ws.Range("F11").Value = 2
$ws.Range("G11").Value = 0.9

# Row 12: 2.2 Design prototypes
$ws.Range("E12").Value = 6
$ws.Range("G12").Value = 0.5

# Row 13: 2.3 Architectural design
$ws.Range("E13").Value = 7
$ws.Range("G13").Value = 0.25

# --- Unhide the Actual Start / Actual Duration / % Complete columns ---
$ws.Columns("E:G").Hidden = $false

# --- View changes: zoom level and active selection ---
$ws.Activate() | Out-Null
$excel.ActiveWindow.Zoom = 75
$ws.Range("E10").Select() | Out-Null
